$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append to the dataset (dash construction start)
$newData = @(
    @(45215, "19:28", 76.7, "natura"),
    @(45217, "19:32", 76.7, "natura"),
    @(45218, "21:24", 76.7, "natura"),
    @(45220, "10:25", 76.7, "natura")
)

$startRow = 33
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $entry = $newData[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
